$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the new "Distance threshold (m)" values for rows 5-9 (column H),
# using full-precision literals so the stored double matches exactly.
$ws.Range("H5:H9").NumberFormat = "0.0"

$ws.Range("H5").Value = 21.566510000000001
$ws.Range("H6").Value = 23.861560000000001
$ws.Range("H7").Value = 59.467800000000004
$ws.Range("H8").Value = 66.31362
$ws.Range("H9").Value = 51.103940000000001

# Update the active selection as recorded in the workbook view.
$ws.Range("J14").Select()
